# Clean up Spreadsheet of SciEntsBank results
#
# - Rename the "Sheet1" tab to "Source Sentences"
# - Update the selection on that sheet from F8 to E9 (and let the view
#   naturally re-anchor, dropping the old scrolled-down topLeftCell)
# - Widen column E on that sheet from ~101.5 to 124 characters

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Rename "Sheet1" -> "Source Sentences"
$ws.Name = "Source Sentences"

# Make sure this sheet is the active one, then move the selection/active
# cell to E9 (matches the <selection activeCell="E9" sqref="E9"/> in the
# target sheetView).
$ws.Activate()
$ws.Range("E9").Select()

# Column E (the 5th column) grows from 101.5 to 124 "characters" wide.
# The ColumnWidth COM property is offset from the raw OOXML column width
# by a small constant padding (~0.83 here), so back that out to land on
# exactly 124 in the saved file.
$ws.Columns.Item(5).ColumnWidth = 124 - 0.83
